# Users can now import sheets of contacts: the "adw3" placeholder in the
# emailAddress column (D4) becomes a real, clickable email address with a
# mailto: hyperlink (adw3@gmail.com), matching what Excel/Office does when a
# proper contact e-mail is filled in so it can be used for importing/emailing
# contacts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell text first.
$ws.Range("D4").Value = "adw3@gmail.com"

# Turn it into a real mailto hyperlink (this also creates the "Hyperlink"
# cell style / font and applies it to D4, just like Excel does natively).
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:adw3@gmail.com")
